$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data / mean calculation
$ws.Range("F5").Value = -3
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = -5
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = -2
$ws.Range("F16").Value = -8
$ws.Range("F17").Value = -13
$ws.Range("F22").Value = 3
